$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 3.311947710001144
$ws.Cells.Item(2, 4).Value = 3.985601948137338
$ws.Cells.Item(2, 5).Value = 11.97137246248165
$ws.Cells.Item(2, 6).Value = 19.9036211589651
$ws.Cells.Item(2, 7).Value = 21.36340211124057
$ws.Cells.Item(2, 8).Value = 12.30490395739066
$ws.Cells.Item(2, 9).Value = 16.89346940137616
$ws.Cells.Item(2, 11).Value = 11.03682614254113
$ws.Cells.Item(2, 13).Value = 14.06451656611418
$ws.Cells.Item(2, 14).Value = 16.4954129498928
$ws.Cells.Item(2, 15).Value = 17.752069031066
$ws.Cells.Item(3, 3).Value = 3.161735716029046
$ws.Cells.Item(3, 4).Value = 3.909766878888276
$ws.Cells.Item(3, 5).Value = 11.86847653828877
$ws.Cells.Item(3, 6).Value = 19.88405102104489
$ws.Cells.Item(3, 7).Value = 21.30153028819468
$ws.Cells.Item(3, 8).Value = 12.34328035612758
$ws.Cells.Item(3, 9).Value = 16.91412700660238
$ws.Cells.Item(3, 11).Value = 10.50750684861218
$ws.Cells.Item(3, 13).Value = 13.77509068870351
$ws.Cells.Item(3, 14).Value = 16.52604496048048
$ws.Cells.Item(3, 15).Value = 17.79561707485923
$ws.Cells.Item(4, 3).Value = 3.06501660086435
$ws.Cells.Item(4, 4).Value = 3.86185860416125
$ws.Cells.Item(4, 5).Value = 11.80982372176826
$ws.Cells.Item(4, 6).Value = 19.8784490040686
$ws.Cells.Item(4, 7).Value = 21.27268973671102
$ws.Cells.Item(4, 8).Value = 12.3691360658243
$ws.Cells.Item(4, 9).Value = 16.93147042600986
$ws.Cells.Item(4, 11).Value = 10.16650597519836
$ws.Cells.Item(4, 13).Value = 13.59707194288078
$ws.Cells.Item(4, 14).Value = 16.54653562220664
$ws.Cells.Item(4, 15).Value = 17.82709848607391
$ws.Cells.Item(5, 3).Value = 3.024502271543396
$ws.Cells.Item(5, 4).Value = 3.842011553345695
$ws.Cells.Item(5, 5).Value = 11.7870853406144
$ws.Cells.Item(5, 6).Value = 19.87778066078012
$ws.Cells.Item(5, 7).Value = 21.26324474823499
$ws.Cells.Item(5, 8).Value = 12.38024823264832
$ws.Cells.Item(5, 9).Value = 16.93970737749336
$ws.Cells.Item(5, 11).Value = 10.02362737531758
$ws.Cells.Item(5, 13).Value = 13.5245531912977
$ws.Cells.Item(5, 14).Value = 16.55530948141188
$ws.Cells.Item(5, 15).Value = 17.84111623613425
$ws.Cells.Item(6, 3).Value = 3.017709344896991
$ws.Cells.Item(6, 4).Value = 3.838696823903164
$ws.Cells.Item(6, 5).Value = 11.78338057634843
$ws.Cells.Item(6, 6).Value = 19.87776721284263
$ws.Cells.Item(6, 7).Value = 21.26181593494803
$ws.Cells.Item(6, 8).Value = 12.38212815278123
$ws.Cells.Item(6, 9).Value = 16.94114566510351
$ws.Cells.Item(6, 11).Value = 9.999669277274648
$ws.Cells.Item(6, 13).Value = 13.51251610226709
$ws.Cells.Item(6, 14).Value = 16.55679198726598
$ws.Cells.Item(6, 15).Value = 17.84351555768341
$ws.Cells.Item(7, 3).Value = 3.064474625759117
$ws.Cells.Item(7, 4).Value = 3.861592232616982
$ws.Cells.Item(7, 5).Value = 11.80951232356228
$ws.Cells.Item(7, 6).Value = 19.87843345244114
$ws.Cells.Item(7, 7).Value = 21.27255300754054
$ws.Cells.Item(7, 8).Value = 12.36928359813058
$ws.Cells.Item(7, 9).Value = 16.93157678118742
$ws.Cells.Item(7, 11).Value = 10.16459477651205
$ws.Cells.Item(7, 13).Value = 13.59609368885531
$ws.Cells.Item(7, 14).Value = 16.54665223290853
$ws.Cells.Item(7, 15).Value = 17.82728272612241
$ws.Cells.Item(8, 3).Value = 3.261101635121028
$ws.Cells.Item(8, 4).Value = 3.959741537279538
$ws.Cells.Item(8, 5).Value = 11.93497056039227
$ws.Cells.Item(8, 6).Value = 19.89554358460099
$ws.Cells.Item(8, 7).Value = 21.34017528521664
$ws.Cells.Item(8, 8).Value = 12.31765993020348
$ws.Cells.Item(8, 9).Value = 16.8996238229809
$ws.Cells.Item(8, 11).Value = 10.85769413403618
$ws.Cells.Item(8, 13).Value = 13.9648480768614
$ws.Cells.Item(8, 14).Value = 16.50562607099579
$ws.Cells.Item(8, 15).Value = 17.76609786588888
$ws.Cells.Item(9, 3).Value = 3.610059325002408
$ws.Cells.Item(9, 4).Value = 4.14095034958194
$ws.Cells.Item(9, 5).Value = 12.21566484253957
$ws.Cells.Item(9, 6).Value = 19.97984556901172
$ws.Cells.Item(9, 7).Value = 21.54487936051726
$ws.Cells.Item(9, 8).Value = 12.23464610655935
$ws.Cells.Item(9, 9).Value = 16.87401536892924
$ws.Cells.Item(9, 11).Value = 12.08620575531405
$ws.Cells.Item(9, 13).Value = 14.68111784484156
$ws.Cells.Item(9, 14).Value = 16.43849438157597
$ws.Cells.Item(9, 15).Value = 17.68390124020691
$ws.Cells.Item(10, 3).Value = 3.843017134189893
$ws.Cells.Item(10, 4).Value = 4.266485686037353
$ws.Cells.Item(10, 5).Value = 12.44119628317852
$ws.Cells.Item(10, 6).Value = 20.07242160860987
$ws.Cells.Item(10, 7).Value = 21.73830560206189
$ws.Cells.Item(10, 8).Value = 12.18480383609528
$ws.Cells.Item(10, 9).Value = 16.87786131249775
$ws.Cells.Item(10, 11).Value = 12.90517442166819
$ws.Cells.Item(10, 13).Value = 15.19754898452706
$ws.Cells.Item(10, 14).Value = 16.3972544055536
$ws.Cells.Item(10, 15).Value = 17.64673841426034
$ws.Cells.Item(11, 3).Value = 3.943742555938747
$ws.Cells.Item(11, 4).Value = 4.321804925726212
$ws.Cells.Item(11, 5).Value = 12.54755678558145
$ws.Cells.Item(11, 6).Value = 20.121091956124
$ws.Cells.Item(11, 7).Value = 21.83537114551681
$ws.Cells.Item(11, 8).Value = 12.16455813639563
$ws.Cells.Item(11, 9).Value = 16.8845340671769
$ws.Cells.Item(11, 11).Value = 13.25899042141726
$ws.Cells.Item(11, 13).Value = 15.42925107631842
$ws.Cells.Item(11, 14).Value = 16.38024002337671
$ws.Cells.Item(11, 15).Value = 17.63490863854826
$ws.Cells.Item(12, 3).Value = 3.981119245729581
$ws.Cells.Item(12, 4).Value = 4.342485127814304
$ws.Cells.Item(12, 5).Value = 12.58833583035586
$ws.Cells.Item(12, 6).Value = 20.14045419732519
$ws.Cells.Item(12, 7).Value = 21.87340424406257
$ws.Cells.Item(12, 8).Value = 12.15724136680099
$ws.Cells.Item(12, 9).Value = 16.8877678867783
$ws.Cells.Item(12, 11).Value = 13.39023939508946
$ws.Cells.Item(12, 13).Value = 15.51644152379235
$ws.Cells.Item(12, 14).Value = 16.37404753157033
$ws.Cells.Item(12, 15).Value = 17.63116090587273
$ws.Cells.Item(13, 3).Value = 3.97310372281714
$ws.Cells.Item(13, 4).Value = 4.338043361826319
$ws.Cells.Item(13, 5).Value = 12.57953160729677
$ws.Cells.Item(13, 6).Value = 20.13624294149744
$ws.Cells.Item(13, 7).Value = 21.86515685104586
$ws.Cells.Item(13, 8).Value = 12.15880159382166
$ws.Cells.Item(13, 9).Value = 16.88704000472761
$ws.Cells.Item(13, 11).Value = 13.36209469143856
$ws.Cells.Item(13, 13).Value = 15.49768928448609
$ws.Cells.Item(13, 14).Value = 16.37537006413404
$ws.Cells.Item(13, 15).Value = 17.63193546563335
$ws.Cells.Item(14, 3).Value = 3.946832963598651
$ws.Cells.Item(14, 4).Value = 4.323511729476976
$ws.Cells.Item(14, 5).Value = 12.55090185452677
$ws.Cells.Item(14, 6).Value = 20.12266629294468
$ws.Cells.Item(14, 7).Value = 21.83847474322414
$ws.Cells.Item(14, 8).Value = 12.1639491660871
$ws.Cells.Item(14, 9).Value = 16.8847859542032
$ws.Cells.Item(14, 11).Value = 13.26984333227962
$ws.Cells.Item(14, 13).Value = 15.43643578727719
$ws.Cells.Item(14, 14).Value = 16.37972554650873
$ws.Cells.Item(14, 15).Value = 17.63458563082918
$ws.Cells.Item(15, 3).Value = 3.930641335527314
$ws.Cells.Item(15, 4).Value = 4.314575472035974
$ws.Cells.Item(15, 5).Value = 12.5334295925207
$ws.Cells.Item(15, 6).Value = 20.11447120102868
$ws.Cells.Item(15, 7).Value = 21.82229651789419
$ws.Cells.Item(15, 8).Value = 12.16714778271189
$ws.Cells.Item(15, 9).Value = 16.88349731432151
$ws.Cells.Item(15, 11).Value = 13.21297970736351
$ws.Cells.Item(15, 13).Value = 15.3988421407503
$ws.Cells.Item(15, 14).Value = 16.38242600892346
$ws.Cells.Item(15, 15).Value = 17.63630430802971
$ws.Cells.Item(16, 3).Value = 3.83632787988635
$ws.Cells.Item(16, 4).Value = 4.26283352143452
$ws.Cells.Item(16, 5).Value = 12.43431784738554
$ws.Cells.Item(16, 6).Value = 20.06937192262865
$ws.Cells.Item(16, 7).Value = 21.73214239049281
$ws.Cells.Item(16, 8).Value = 12.18617585636174
$ws.Cells.Item(16, 9).Value = 16.87752423340931
$ws.Cells.Item(16, 11).Value = 12.88167115045331
$ws.Cells.Item(16, 13).Value = 15.18233430807751
$ws.Cells.Item(16, 14).Value = 16.39840141829215
$ws.Cells.Item(16, 15).Value = 17.64761384140504
$ws.Cells.Item(17, 3).Value = 3.777116790482404
$ws.Cells.Item(17, 4).Value = 4.230625665571338
$ws.Cells.Item(17, 5).Value = 12.37445310327328
$ws.Cells.Item(17, 6).Value = 20.04337703348558
$ws.Cells.Item(17, 7).Value = 21.67914083277099
$ws.Cells.Item(17, 8).Value = 12.19847130701029
$ws.Cells.Item(17, 9).Value = 16.87512042773935
$ws.Cells.Item(17, 11).Value = 12.67359466087782
$ws.Cells.Item(17, 13).Value = 15.04862462333775
$ws.Cells.Item(17, 14).Value = 16.40864857117289
$ws.Cells.Item(17, 15).Value = 17.65585350523357
$ws.Cells.Item(18, 3).Value = 3.742566663213184
$ws.Cells.Item(18, 4).Value = 4.211932925443588
$ws.Cells.Item(18, 5).Value = 12.34037728031914
$ws.Cells.Item(18, 6).Value = 20.02904340552377
$ws.Cells.Item(18, 7).Value = 21.64951183521059
$ws.Cells.Item(18, 8).Value = 12.20577180000612
$ws.Cells.Item(18, 9).Value = 16.87420131903283
$ws.Cells.Item(18, 11).Value = 12.55215273368986
$ws.Cells.Item(18, 13).Value = 14.97141962482674
$ws.Cells.Item(18, 14).Value = 16.41470682780427
$ws.Cells.Item(18, 15).Value = 17.66107040298369
$ws.Cells.Item(19, 3).Value = 3.730784204712092
$ws.Cells.Item(19, 4).Value = 4.205575440319516
$ws.Cells.Item(19, 5).Value = 12.32890223547026
$ws.Cells.Item(19, 6).Value = 20.02429672870339
$ws.Cells.Item(19, 7).Value = 21.63962781626586
$ws.Cells.Item(19, 8).Value = 12.20828283966755
$ws.Cells.Item(19, 9).Value = 16.87396974877523
$ws.Cells.Item(19, 11).Value = 12.51073331238121
$ws.Cells.Item(19, 13).Value = 14.94523075160558
$ws.Cells.Item(19, 14).Value = 16.41678629677934
$ws.Cells.Item(19, 15).Value = 17.66291872383389
$ws.Cells.Item(20, 3).Value = 3.783471072413944
$ws.Cells.Item(20, 4).Value = 4.234071678079117
$ws.Cells.Item(20, 5).Value = 12.38078915811899
$ws.Cells.Item(20, 6).Value = 20.04608035387253
$ws.Cells.Item(20, 7).Value = 21.68469453244336
$ws.Cells.Item(20, 8).Value = 12.19713878363122
$ws.Cells.Item(20, 9).Value = 16.87532835157185
$ws.Cells.Item(20, 11).Value = 12.69592741026793
$ws.Cells.Item(20, 13).Value = 15.06288980694916
$ws.Cells.Item(20, 14).Value = 16.40754073720423
$ws.Cells.Item(20, 15).Value = 17.65492692651999
$ws.Cells.Item(21, 3).Value = 3.954570185402754
$ws.Cells.Item(21, 4).Value = 4.327787379011769
$ws.Cells.Item(21, 5).Value = 12.55929777991842
$ws.Cells.Item(21, 6).Value = 20.12662888994033
$ws.Cells.Item(21, 7).Value = 21.84627752441219
$ws.Cells.Item(21, 8).Value = 12.16242769927006
$ws.Cells.Item(21, 9).Value = 16.88542884694081
$ws.Cells.Item(21, 11).Value = 13.29701425085376
$ws.Cells.Item(21, 13).Value = 15.45444301535111
$ws.Cells.Item(21, 14).Value = 16.37843944219085
$ws.Cells.Item(21, 15).Value = 17.6337873340673
$ws.Cells.Item(22, 3).Value = 4.061925705840778
$ws.Cells.Item(22, 4).Value = 4.387469386375469
$ws.Cells.Item(22, 5).Value = 12.67887455045753
$ws.Cells.Item(22, 6).Value = 20.1846972340127
$ws.Cells.Item(22, 7).Value = 21.95930627402494
$ws.Cells.Item(22, 8).Value = 12.1417814279619
$ws.Cells.Item(22, 9).Value = 16.89615006952428
$ws.Cells.Item(22, 11).Value = 13.67391368289057
$ws.Cells.Item(22, 13).Value = 15.70710301354025
$ws.Cells.Item(22, 14).Value = 16.36087982163114
$ws.Cells.Item(22, 15).Value = 17.62423837094299
$ws.Cells.Item(23, 3).Value = 4.005039980012733
$ws.Cells.Item(23, 4).Value = 4.355762747493383
$ws.Cells.Item(23, 5).Value = 12.61480088066889
$ws.Cells.Item(23, 6).Value = 20.15321271893874
$ws.Cells.Item(23, 7).Value = 21.898311541027
$ws.Cells.Item(23, 8).Value = 12.15261389467937
$ws.Cells.Item(23, 9).Value = 16.89005146745344
$ws.Cells.Item(23, 11).Value = 13.47422532316333
$ws.Cells.Item(23, 13).Value = 15.57257744818677
$ws.Cells.Item(23, 14).Value = 16.37011834161851
$ws.Cells.Item(23, 15).Value = 17.62894383934093
$ws.Cells.Item(24, 3).Value = 3.780599885755631
$ws.Cells.Item(24, 4).Value = 4.232514283230761
$ws.Cells.Item(24, 5).Value = 12.37792356193438
$ws.Cells.Item(24, 6).Value = 20.04485627812049
$ws.Cells.Item(24, 7).Value = 21.68218107897674
$ws.Cells.Item(24, 8).Value = 12.19774049584379
$ws.Cells.Item(24, 9).Value = 16.87523290732557
$ws.Cells.Item(24, 11).Value = 12.68583642763951
$ws.Cells.Item(24, 13).Value = 15.05644154964117
$ws.Cells.Item(24, 14).Value = 16.40804106859579
$ws.Cells.Item(24, 15).Value = 17.65534433834379
$ws.Cells.Item(25, 3).Value = 3.519703995121711
$ws.Cells.Item(25, 4).Value = 4.093212861306198
$ws.Cells.Item(25, 5).Value = 12.13619838315561
$ws.Cells.Item(25, 6).Value = 19.95163090620493
$ws.Cells.Item(25, 7).Value = 21.48186367391427
$ws.Cells.Item(25, 8).Value = 12.25514865932113
$ws.Cells.Item(25, 9).Value = 16.87696651526488
$ws.Cells.Item(25, 11).Value = 11.7683076412154
$ws.Cells.Item(25, 13).Value = 14.48867973737988
$ws.Cells.Item(25, 14).Value = 16.45523326560823
$ws.Cells.Item(25, 15).Value = 17.70207101633056